# Swap the presentation's applied design theme from the custom "Integral"
# theme (ppt/theme/theme1.xml) to the built-in "Office Theme" palette, i.e.
# re-colour the 12 theme colour slots (dk1, lt1, dk2, lt2, accent1-6, hlink,
# folHlink) of the active slide master's theme to match the stock Office
# theme colour scheme.
#
# PowerPoint stores RGB() values in COM as 0x00BBGGRR, so convert each
# target hex colour (RRGGBB) into that packed integer before assigning it.

function Convert-HexToRgbInt($hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    return ($b * 65536) + ($g * 256) + $r
}

# Office Theme colour scheme, in a:clrScheme child order.
$officeThemeColors = @(
    "000000", # dk1
    "FFFFFF", # lt1
    "44546A", # dk2
    "E7E6E6", # lt2
    "5B9BD5", # accent1
    "ED7D31", # accent2
    "A5A5A5", # accent3
    "FFC000", # accent4
    "4472C4", # accent5
    "70AD47", # accent6
    "0563C1", # hlink
    "954F72"  # folHlink
)

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$themeColors = $master.Theme.ThemeColorScheme

for ($i = 1; $i -le $officeThemeColors.Count; $i++) {
    $themeColors.Item($i).RGB = Convert-HexToRgbInt $officeThemeColors[$i - 1]
}
